$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert the two new rows that hold the new error-code entries.
#    - row 4  : new "#1002 / Erreur lors de la lecture de la ressource."
#               entry, extending the "Generique" (#1) merged group.
#    - row 8  : new "#3003 / Impossible de recuperer les donnees du test."
#               entry, extending the "Application" (#3) merged group.
# ------------------------------------------------------------------
$ws.Rows(4).EntireRow.Insert()
$ws.Rows(8).EntireRow.Insert()

# ------------------------------------------------------------------
# 2. Fix up the borders/styles for the rows that moved around.
#    Row 3 used to be the *bottom* of the A2:A3 merge (border on
#    bottom only); now that the merge grows to A2:A4 it becomes the
#    *middle* row, so it needs side-only borders. Row 4 becomes the
#    new bottom of that merge, so it gets the old "bottom border"
#    look that row 3 used to have.
# ------------------------------------------------------------------
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)

$ws.Range("A3").Borders.Item(8).LineStyle = 0
$ws.Range("A3").Borders.Item(9).LineStyle = 0
$ws.Range("A3").Borders.Item(7).LineStyle = 1
$ws.Range("A3").Borders.Item(10).LineStyle = 1

# Row 7 (old row 6, the middle of the Application group) keeps its
# normal box-border look; row 8 (brand new) should look the same.
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)

$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial(-4122)

$ws.Range("C7").Copy()
$ws.Range("C8").PasteSpecial(-4122)

$ws.Range("B4").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C4").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 3. Cell values.
# ------------------------------------------------------------------

# "Generique" category header (was "Non trouve")
$ws.Range("A2").Value = "Générique" + [char]10 + "#1"

# Row 3 unchanged content (#1001 / Fichier non trouvé.)
$ws.Range("B3").Value = "#1001"
$ws.Range("C3").Value = "Fichier non trouvé."

# Row 4 - brand new entry
$ws.Range("B4").Value = "#1002"
$ws.Range("C4").Value = "Erreur lors de la lecture de la ressource."

# Row 5 (was row 4) unchanged
$ws.Range("A5").Value = "Versions" + [char]10 + "#2"
$ws.Range("B5").Value = "#2001"
$ws.Range("C5").Value = "Historique des versions introuvable."

# Row 6 (was row 5) unchanged
$ws.Range("A6").Value = "Application" + [char]10 + "#3"
$ws.Range("B6").Value = "#3001"
$ws.Range("C6").Value = "Application inconnue."

# Row 7 (was row 6) unchanged
$ws.Range("B7").Value = "#3002"
$ws.Range("C7").Value = "Contenu de l'application inacessible."

# Row 8 - brand new entry
$ws.Range("B8").Value = "#3003"
$ws.Range("C8").Value = "Impossible de récupérer les données du test."

# Row 9 (was row 7) - code changes from #3003 to #3004, message unchanged
$ws.Range("B9").Value = "#3004"
$ws.Range("C9").Value = "Impossible de récupérer la correction."

# ------------------------------------------------------------------
# 4. Merged cells for the two categories.
# ------------------------------------------------------------------
$ws.Range("A2:A3").UnMerge()
$ws.Range("A2:A4").Merge()

$ws.Range("A6:A9").UnMerge()
$ws.Range("A6:A9").Merge()

# ------------------------------------------------------------------
# 5. Selection / view bits.
# ------------------------------------------------------------------
$ws.Range("C4").Select()

$aw = $excel.ActiveWindow
$aw.WindowState = -4143
$aw.Width = 16200
$aw.Height = 9360
